# Insert a new data row before row 39 (pushes existing rows 39..122 down to 40..123)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(39).Insert()

$ws.Range("A39").Value = 8
$ws.Range("B39").Value = "Terminal La Palmera de La Serena"
$ws.Range("C39").Value = "Coquimbo"
$ws.Range("D39").Value = 44987
$ws.Range("E39").Value = 4
$ws.Range("F39").Value = 100112030
$ws.Range("G39").Value = "Poroto granado"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 400
$ws.Range("K39").Value = 36000
$ws.Range("L39").Value = 37000
$ws.Range("M39").Value = 36500
$ws.Range("N39").Value = "`$/malla 25 kilos"
$ws.Range("O39").Value = "Provincia del Elquí"
$ws.Range("P39").Value = 1460
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"
